$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as text even when the value looks
# like a plain number (e.g. "608.92"), matching the original workbook's
# inline-string cryptocurrency price formatting, and then strip the
# temporary text number-format back off so no stray style is left behind.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.791.61"
$ws.Range("E2").Value = "  -3.43%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.141.01"
$ws.Range("E3").Value = "  -3.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "608.92"
$ws.Range("E5").Value = "  +0.53%  "

# Row 6 - Solana
Set-TextValue "D6" "145.73"
$ws.Range("E6").Value = "  -7.06%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.134.55"
$ws.Range("E8").Value = "  -3.36%  "

# Row 9 - XRP
Set-TextValue "D9" "0.522"
$ws.Range("E9").Value = "  -2.85%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.151"
$ws.Range("E10").Value = "  -5.56%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.29"
$ws.Range("E11").Value = "  -6.77%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.472"
$ws.Range("E12").Value = "  -3.45%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000252"
$ws.Range("E13").Value = "  -4.95%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.41"
$ws.Range("E14").Value = "  -7.51%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.643.28"
$ws.Range("E15").Value = "  -3.74%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.828.26"
$ws.Range("E17").Value = "  -3.53%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.134.45"
$ws.Range("E18").Value = "  -3.57%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.83"
$ws.Range("E19").Value = "  -5.27%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "477.52"
$ws.Range("E20").Value = "  -3.80%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.62"
$ws.Range("E21").Value = "  -4.06%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.708"
$ws.Range("E22").Value = "  -4.40%  "

# Row 23 - Uniswap
Set-TextValue "D23" "7.81"
$ws.Range("E23").Value = "  -2.04%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.54"
$ws.Range("E24").Value = "  -6.40%  "

# Row 25 - Litecoin
Set-TextValue "D25" "83.72"
$ws.Range("E25").Value = "  -3.32%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.14%  "

# Row 27 - PancakeSwap
Set-TextValue "D27" "2.79"
$ws.Range("E27").Value = "  -7.21%  "

# Row 28 - RenderToken
Set-TextValue "D28" "8.47"
$ws.Range("E28").Value = "  -6.10%  "

# Row 29 - Hedera
Set-TextValue "D29" "0.120"
$ws.Range("E29").Value = "  -9.57%  "

# Row 30 - NEARProtocol
Set-TextValue "D30" "6.88"
$ws.Range("E30").Value = "  -0.95%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "2.08"
$ws.Range("E31").Value = "  -11.60%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  +0.01%  "

# Row 33 - Stacks
Set-TextValue "D33" "2.69"
$ws.Range("E33").Value = "  -4.64%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "26.25"
$ws.Range("E34").Value = "  -5.07%  "

# Row 35 - Mantle
Set-TextValue "D35" "1.12"
$ws.Range("E35").Value = "  -0.76%  "

# Rows 36/37 - PEPE and Filecoin swap order (PEPE now listed first, then
# Filecoin), each with refreshed data.
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0791"
$ws.Range("E36").Value = "  +2.94%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D37" "5.94"
$ws.Range("E37").Value = "  -6.48%  "

# Row 38 - OKB
Set-TextValue "D38" "52.91"
$ws.Range("E38").Value = "  -4.99%  "

# Row 39 - Bittensor
Set-TextValue "D39" "457.09"
$ws.Range("E39").Value = "  -6.59%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "2.99"
$ws.Range("E40").Value = "  -11.73%  "

# Row 41 - VeChain
Set-TextValue "D41" "0.0395"
$ws.Range("E41").Value = "  -5.66%  "

# Row 42 - Kaspa
Set-TextValue "D42" "0.118"
$ws.Range("E42").Value = "  -8.32%  "

# Row 43 - Cosmos
Set-TextValue "D43" "8.34"
$ws.Range("E43").Value = "  -4.09%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.844.43"
$ws.Range("E44").Value = "  -4.75%  "

# Row 45 - Fetch.AI
Set-TextValue "D45" "2.29"
$ws.Range("E45").Value = "  -9.26%  "

# Row 46 - TheGraph
Set-TextValue "D46" "0.266"
$ws.Range("E46").Value = "  -7.70%  "

# Row 47 - ThetaToken
Set-TextValue "D47" "2.44"
$ws.Range("E47").Value = "  +0.28%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "26.25"
$ws.Range("E49").Value = "  -7.04%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -3.83%  "

# Row 51 - Monero
Set-TextValue "D51" "119.18"
$ws.Range("E51").Value = "  -2.33%  "
